$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "avatar" column entirely (E)
$ws.Columns("E").Delete()

# Update row 2 with the new student record
$ws.Range("A2").Value = "HE838383"
$ws.Range("B2").Value = "Viruss"
$ws.Range("C2").Value = "AI1908"
$ws.Range("D2").Value = "AI"

# Remove row 3 entirely (the old "HE123457 / trump" record)
$ws.Rows("3").Delete()
